$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'36.414.64"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = "'  +0.19%  "
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.Value = "'1.954.69"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.Value = "'  -1.16%  "
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.Value = "'  -0.04%  "
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.Value = "'243.83"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = "'  -0.24%  "
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.Value = "'0.617"
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.Value = "'  -1.39%  "
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.Value = "'58.07"
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.Value = "'  -3.60%  "
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.Value = "'  -0.03%  "
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.Value = "'0.365"
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = "'  -2.64%  "
$c.Style = "Normal"
$c = $ws.Range("B10")
$c.Value = "'Dogecoin"
$c.Style = "Normal"
$c = $ws.Range("C10")
$c.Value = "'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.Value = "'0.0855"
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.Value = "'  +6.16%  "
$c.Style = "Normal"
$c = $ws.Range("B11")
$c.Value = "'TRON"
$c.Style = "Normal"
$c = $ws.Range("C11")
$c.Value = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.Value = "'0.104"
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.Value = "'  +0.41%  "
$c.Style = "Normal"
$c = $ws.Range("B12")
$c.Value = "'WrappedliquidstakedEther2.0"
$c.Style = "Normal"
$c = $ws.Range("C12")
$c.Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.Value = "'2.241.11"
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.Value = "'  -1.09%  "
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.Value = "'21.62"
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.Value = "'  -5.62%  "
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.Value = "'0.824"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.Value = "'  -3.96%  "
$c.Style = "Normal"
$c = $ws.Range("B15")
$c.Value = "'Chainlink"
$c.Style = "Normal"
$c = $ws.Range("C15")
$c.Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.Value = "'13.56"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.Value = "'  -4.09%  "
$c.Style = "Normal"
$c = $ws.Range("B16")
$c.Value = "'Polkadot"
$c.Style = "Normal"
$c = $ws.Range("C16")
$c.Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$c.Style = "Normal"
$c = $ws.Range("D16")
$c.Value = "'5.21"
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.Value = "'  -4.14%  "
$c.Style = "Normal"
$c = $ws.Range("B17")
$c.Value = "'WrappedEther"
$c.Style = "Normal"
$c = $ws.Range("C17")
$c.Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.Value = "'1.945.43"
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.Value = "'  -1.15%  "
$c.Style = "Normal"
$c = $ws.Range("B18")
$c.Value = "'WrappedBTC"
$c.Style = "Normal"
$c = $ws.Range("C18")
$c.Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.Value = "'36.373.62"
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.Value = "'  +0.17%  "
$c.Style = "Normal"
$c = $ws.Range("B19")
$c.Value = "'Litecoin"
$c.Style = "Normal"
$c = $ws.Range("C19")
$c.Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.Value = "'70.05"
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.Value = "'  -1.82%  "
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.Value = "'0.0₃0875"
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.Value = "'  +1.70%  "
$c.Style = "Normal"
$c = $ws.Range("B21")
$c.Value = "'BitcoinCash"
$c.Style = "Normal"
$c = $ws.Range("C21")
$c.Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.Value = "'229.76"
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.Value = "'  -2.19%  "
$c.Style = "Normal"
$c = $ws.Range("B22")
$c.Value = "'Uniswap"
$c.Style = "Normal"
$c = $ws.Range("C22")
$c.Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.Value = "'5.06"
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.Value = "'  -3.90%  "
$c.Style = "Normal"
$c = $ws.Range("B23")
$c.Value = "'Dai"
$c.Style = "Normal"
$c = $ws.Range("C23")
$c.Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.Value = "'1.00"
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.Value = "'  -0.09%  "
$c.Style = "Normal"
$c = $ws.Range("B24")
$c.Value = "'PancakeSwap"
$c.Style = "Normal"
$c = $ws.Range("C24")
$c.Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.Value = "'2.44"
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.Value = "'  -5.24%  "
$c.Style = "Normal"
$c = $ws.Range("B25")
$c.Value = "'Toncoin"
$c.Style = "Normal"
$c = $ws.Range("C25")
$c.Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.Value = "'2.31"
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.Value = "'  +0.55%  "
$c.Style = "Normal"
$c = $ws.Range("B26")
$c.Value = "'Cosmos"
$c.Style = "Normal"
$c = $ws.Range("C26")
$c.Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c.Style = "Normal"
$c = $ws.Range("D26")
$c.Value = "'9.32"
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.Value = "'  -4.86%  "
$c.Style = "Normal"
$c = $ws.Range("B27")
$c.Value = "'Monero"
$c.Style = "Normal"
$c = $ws.Range("C27")
$c.Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.Value = "'161.52"
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.Value = "'  +0.02%  "
$c.Style = "Normal"
$c = $ws.Range("B28")
$c.Value = "'Kaspa"
$c.Style = "Normal"
$c = $ws.Range("C28")
$c.Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c.Style = "Normal"
$c = $ws.Range("D28")
$c.Value = "'0.132"
$c.Style = "Normal"
$c = $ws.Range("E28")
$c.Value = "'  -2.43%  "
$c.Style = "Normal"
$c = $ws.Range("B29")
$c.Value = "'EthereumClassic"
$c.Style = "Normal"
$c = $ws.Range("C29")
$c.Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c.Style = "Normal"
$c = $ws.Range("D29")
$c.Value = "'19.50"
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.Value = "'  -1.70%  "
$c.Style = "Normal"
$c = $ws.Range("B30")
$c.Value = "'Stellar"
$c.Style = "Normal"
$c = $ws.Range("C30")
$c.Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.Value = "'0.118"
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.Value = "'  -1.46%  "
$c.Style = "Normal"
$c = $ws.Range("B31")
$c.Value = "'ImmutableX"
$c.Style = "Normal"
$c = $ws.Range("C31")
$c.Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.Value = "'1.17"
$c.Style = "Normal"
$c = $ws.Range("E31")
$c.Value = "'  +2.12%  "
$c.Style = "Normal"
$c = $ws.Range("B32")
$c.Value = "'Filecoin"
$c.Style = "Normal"
$c = $ws.Range("C32")
$c.Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c.Style = "Normal"
$c = $ws.Range("D32")
$c.Value = "'4.68"
$c.Style = "Normal"
$c = $ws.Range("E32")
$c.Value = "'  -4.41%  "
$c.Style = "Normal"
$c = $ws.Range("B33")
$c.Value = "'Hedera"
$c.Style = "Normal"
$c = $ws.Range("C33")
$c.Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c.Style = "Normal"
$c = $ws.Range("D33")
$c.Value = "'0.0647"
$c.Style = "Normal"
$c = $ws.Range("E33")
$c.Value = "'  +4.03%  "
$c.Style = "Normal"
$c = $ws.Range("B34")
$c.Value = "'InternetComputer(DFINITY)"
$c.Style = "Normal"
$c = $ws.Range("C34")
$c.Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c.Style = "Normal"
$c = $ws.Range("D34")
$c.Value = "'4.30"
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.Value = "'  -4.08%  "
$c.Style = "Normal"
$c = $ws.Range("B35")
$c.Value = "'THORChain"
$c.Style = "Normal"
$c = $ws.Range("C35")
$c.Value = "'https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$c.Style = "Normal"
$c = $ws.Range("D35")
$c.Value = "'6.32"
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.Value = "'  +3.60%  "
$c.Style = "Normal"
$c = $ws.Range("B36")
$c.Value = "'BinanceUSD"
$c.Style = "Normal"
$c = $ws.Range("C36")
$c.Value = "'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$c.Style = "Normal"
$c = $ws.Range("D36")
$c.Value = "'1.00"
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.Value = "'  -0.05%  "
$c.Style = "Normal"
$c = $ws.Range("B37")
$c.Value = "'WEMIXToken"
$c.Style = "Normal"
$c = $ws.Range("C37")
$c.Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$c.Style = "Normal"
$c = $ws.Range("D37")
$c.Value = "'1.78"
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.Value = "'  -0.45%  "
$c.Style = "Normal"
$c = $ws.Range("B38")
$c.Value = "'LidoDAOToken"
$c.Style = "Normal"
$c = $ws.Range("C38")
$c.Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$c.Style = "Normal"
$c = $ws.Range("D38")
$c.Value = "'2.17"
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.Value = "'  -5.17%  "
$c.Style = "Normal"
$c = $ws.Range("B39")
$c.Value = "'RenderToken"
$c.Style = "Normal"
$c = $ws.Range("C39")
$c.Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.Value = "'3.04"
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.Value = "'  -0.77%  "
$c.Style = "Normal"
$c = $ws.Range("B40")
$c.Value = "'Cronos"
$c.Style = "Normal"
$c = $ws.Range("C40")
$c.Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.Value = "'0.0979"
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.Value = "'  +0.00%  "
$c.Style = "Normal"
$c = $ws.Range("B41")
$c.Value = "'HuobiToken"
$c.Style = "Normal"
$c = $ws.Range("C41")
$c.Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.Value = "'2.89"
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.Value = "'  -0.04%  "
$c.Style = "Normal"
$c = $ws.Range("B42")
$c.Value = "'TrustWalletToken"
$c.Style = "Normal"
$c = $ws.Range("C42")
$c.Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.Value = "'1.17"
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.Value = "'  -5.33%  "
$c.Style = "Normal"
$c = $ws.Range("B43")
$c.Value = "'VeChain"
$c.Style = "Normal"
$c = $ws.Range("C43")
$c.Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.Value = "'0.0211"
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.Value = "'  -1.11%  "
$c.Style = "Normal"
$c = $ws.Range("B44")
$c.Value = "'InjectiveProtocol"
$c.Style = "Normal"
$c = $ws.Range("C44")
$c.Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.Value = "'15.89"
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.Value = "'  -1.91%  "
$c.Style = "Normal"
$c = $ws.Range("B45")
$c.Value = "'Maker"
$c.Style = "Normal"
$c = $ws.Range("C45")
$c.Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.Value = "'1.358.66"
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.Value = "'  -0.61%  "
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.Value = "'1.03"
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.Value = "'  -5.74%  "
$c.Style = "Normal"
$c = $ws.Range("B47")
$c.Value = "'Aave"
$c.Style = "Normal"
$c = $ws.Range("C47")
$c.Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.Value = "'88.07"
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.Value = "'  -4.86%  "
$c.Style = "Normal"
$c = $ws.Range("B48")
$c.Value = "'FraxShare"
$c.Style = "Normal"
$c = $ws.Range("C48")
$c.Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.Value = "'7.16"
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.Value = "'  -5.71%  "
$c.Style = "Normal"
$c = $ws.Range("B49")
$c.Value = "'MXToken"
$c.Style = "Normal"
$c = $ws.Range("C49")
$c.Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.Value = "'2.82"
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.Value = "'  -0.86%  "
$c.Style = "Normal"
$c = $ws.Range("B50")
$c.Value = "'MultiversX"
$c.Style = "Normal"
$c = $ws.Range("C50")
$c.Value = "'https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.Value = "'45.41"
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.Value = "'  +2.07%  "
$c.Style = "Normal"
$c = $ws.Range("B51")
$c.Value = "'RocketPoolETH"
$c.Style = "Normal"
$c = $ws.Range("C51")
$c.Value = "'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.Value = "'2.131.73"
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.Value = "'  -1.27%  "
$c.Style = "Normal"
